$d = $word.ActiveDocument

# 1. Merge the split "doorgevoer" / _GoBack bookmark / "d." runs in the
#    "Gebruikersnaam" userstory into a single run "...doorgevoerd." and
#    drop the bookmark. Find/Replace (wildcards on) merges runs across the
#    bookmark split and removes the bookmark automatically.
$old1 = "Als beheerder wil ik dat na het wijzigen van de Gebruikersnaam, deze wijziging direct zichtbaar is, zodat de gebruiker weet dat de wijziging goed is doorgevoerd."
$new1 = "Als beheerder wil ik dat na het wijzigen van de Gebruikersnaam, deze wijziging direct zichtbaar is, zodat de gebruiker weet dat de wijziging goed is doorgevoerd."
$d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2) | Out-Null

# 2. Move the "lijst van gebruikers" userstory paragraph: insert it right
#    after the "... rollen aan gebruikers toe kunnen wijzen ..." paragraph
#    (just before the "OPTIONEEL: Meertaligheid..." Kop2 heading), as a
#    single run, with a collapsed _GoBack bookmark at the end of the text.
$newText = "Als beheerder wil ik dat ik in een lijst van gebruikers, de rollen van gebruikers wijzigen, zodat ik een gebruiker meer of minder rechten kan geven."

$anchorText = "Als beheerder wil ik bij het aanmaken van gebruikers, rollen aan gebruikers toe kunnen wijzen, zodat ik kan bepalen wie wat mag doen binnen de applicatie."
$anchor = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.TrimEnd([char]13) -eq $anchorText) {
        $anchor = $d.Paragraphs($i)
        break
    }
}

$anchor.Range.InsertParagraphAfter() | Out-Null
$newParaIndex = $anchor.Index + 1
$newPara = $d.Paragraphs($newParaIndex)
# Append a sentinel character so the bookmark range is non-collapsed while
# we create it (collapsed ranges aren't placed correctly by Bookmarks.Add),
# then strip the sentinel back out so the bookmark ends up zero-width,
# sitting right after the text and before the paragraph mark.
$newPara.Range.Text = $newText + "X"
$sentinelStart = $newPara.Range.End - 2
$sentinelRange = $d.Range($sentinelStart, $sentinelStart + 1)
$d.Bookmarks.Add("_GoBack", $sentinelRange) | Out-Null
$sentinelRange2 = $d.Range($sentinelStart, $sentinelStart + 1)
$sentinelRange2.Text = ""

# 3. Remove the old "lijst van gebruikers" paragraph further down the
#    document (after "... zichzelf kunnen registreren ..."). Search after
#    the paragraph we just inserted, since its text now matches too.
$old = $null
for ($i = $newParaIndex + 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.TrimEnd([char]13) -eq $newText) {
        $old = $d.Paragraphs($i)
        break
    }
}
$old.Range.Delete()
